$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7452.778
$ws.Range("I141").Value = 8009.375
$ws.Range("K141").Value = 24028.125
$ws.Range("M141").Value = -18848.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 16980
$ws.Range("J9").Value = 16980
$ws.Range("L9").Value = 16980
$ws.Range("N9").Value = -17320
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0
$ws.Range("H16").Value = 2779.6667
$ws.Range("I16").Value = 2503
$ws.Range("J16").Value = 3333
$ws.Range("K16").Value = 2503
$ws.Range("L16").Value = 3333
$ws.Range("M16").Value = -2216
$ws.Range("N16").Value = -3907
$ws.Range("H20").Value = 16980
$ws.Range("J20").Value = 16980
$ws.Range("L20").Value = 16980
$ws.Range("N20").Value = -17520
$ws.Range("H22").Value = 5500
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -701
$ws.Range("N22").Value = -10598
$ws.Range("H23").Value = 8740
$ws.Range("J23").Value = 8740
$ws.Range("L23").Value = 8740
$ws.Range("N23").Value = -9258
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -598
$ws.Range("H30").Value = 3196.6667
$ws.Range("I30").Value = 2800
$ws.Range("J30").Value = 3395
$ws.Range("K30").Value = 2800
$ws.Range("L30").Value = 3395
$ws.Range("M30").Value = -2650
$ws.Range("N30").Value = -3695
$ws.Range("H45").Value = 2577.0908
$ws.Range("I45").Value = 2459.4138
$ws.Range("J45").Value = 2804.6
$ws.Range("K45").Value = 2459.4138
$ws.Range("L45").Value = 2804.6
$ws.Range("M45").Value = -2082.4138
$ws.Range("N45").Value = -3558.6
$ws.Range("H46").Value = 4975.3335
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 5963
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 5963
$ws.Range("M46").Value = -2681
$ws.Range("N46").Value = -6601
$ws.Range("H50").Value = 1314.7142
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1314.7142
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").Value = 1314.7142
$ws.Range("N50").Value = -2742.7142
$ws.Range("H102").Value = 3490.7407
$ws.Range("I102").Value = 3324.375
$ws.Range("J102").Value = 3732.7273
$ws.Range("K102").Value = 3324.375
$ws.Range("L102").Value = 3732.7273
$ws.Range("M102").Value = -1702.375
$ws.Range("N102").Value = -6976.7273
$ws.Range("H122").Value = 3038.963
$ws.Range("I122").Value = 2952.2083
$ws.Range("K122").Value = 8856.624899999999
$ws.Range("M122").Value = -6406.624899999999
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3157.6155
$ws.Range("I86").Value = 2928.5
$ws.Range("J86").Value = 3424.9167
$ws.Range("K86").Value = 2928.5
$ws.Range("L86").Value = 3424.9167
$ws.Range("M86").Value = -1805.5
$ws.Range("N86").Value = -5670.9167
$ws.Range("H89").Value = 3157.6155
$ws.Range("I89").Value = 2928.5
$ws.Range("J89").Value = 3424.9167
$ws.Range("K89").Value = 14642.5
$ws.Range("L89").Value = 17124.5835
$ws.Range("M89").Value = -9026.5
$ws.Range("N89").Value = -28356.5835
$ws.Range("H99").Value = 1147.3103
$ws.Range("I99").Value = 842.5
$ws.Range("J99").Value = 1824.6666
$ws.Range("K99").Value = 842.5
$ws.Range("L99").Value = 1824.6666
$ws.Range("M99").Value = 655.5
$ws.Range("N99").Value = -4820.6666
$ws.Range("H105").Value = 1680
$ws.Range("I105").Value = 1643
$ws.Range("J105").Value = 1703.125
$ws.Range("K105").Value = 1643
$ws.Range("L105").Value = 1703.125
$ws.Range("M105").Value = 104
$ws.Range("N105").Value = -5197.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 6398.3335
$ws.Range("J56").Value = 9200
$ws.Range("L56").Value = 9200
$ws.Range("N56").Value = -10890
$ws.Range("H58").Value = 1951.25
$ws.Range("I58").Value = 1725.5
$ws.Range("J58").Value = 2628.5
$ws.Range("K58").Value = 1725.5
$ws.Range("L58").Value = 2628.5
$ws.Range("M58").Value = -1522.5
$ws.Range("N58").Value = -3034.5
$ws.Range("H105").Value = 1197.1428
$ws.Range("I105").Value = 836.6667
$ws.Range("K105").Value = 836.6667
$ws.Range("M105").Value = 910.3333
$ws.Range("H136").Value = 1951.25
$ws.Range("I136").Value = 1725.5
$ws.Range("J136").Value = 2628.5
$ws.Range("K136").Value = 5176.5
$ws.Range("L136").Value = 7885.5
$ws.Range("M136").Value = -2626.5
$ws.Range("N136").Value = -12985.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 839.75
$ws.Range("J131").Value = 835.2841
$ws.Range("L131").Value = 2505.8523
$ws.Range("N131").Value = -12585.8523

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 413.625
$ws.Range("I107").Value = 375.2
$ws.Range("J107").Value = 990
$ws.Range("K107").Value = 375.2
$ws.Range("L107").Value = 990
$ws.Range("M107").Value = 1544.8
$ws.Range("N107").Value = -4830
$ws.Range("H122").Value = 2795.92
$ws.Range("I122").Value = 2946.5881
$ws.Range("J122").Value = 2475.75
$ws.Range("K122").Value = 8839.764299999999
$ws.Range("L122").Value = 7427.25
$ws.Range("M122").Value = -6389.764299999999
$ws.Range("N122").Value = -12327.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3012
$ws.Range("I122").Value = 2499.6
$ws.Range("J122").Value = 3524.4
$ws.Range("K122").Value = 7498.799999999999
$ws.Range("L122").Value = 10573.2
$ws.Range("M122").Value = -5048.799999999999
$ws.Range("N122").Value = -15473.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -844
$ws.Range("H18").Value = 1999
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 1998
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 1998
$ws.Range("M18").Value = -1827
$ws.Range("N18").Value = -2344
$ws.Range("H19").Value = 1974.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1974.5
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = 1974.5
$ws.Range("N19").Value = -2322.5
$ws.Range("H122").Value = 1532.3462
$ws.Range("I122").Value = 1516.3334
$ws.Range("J122").Value = 1554.1818
$ws.Range("K122").Value = 4549.0002
$ws.Range("L122").Value = 4662.5454
$ws.Range("M122").Value = -2099.0002
$ws.Range("N122").Value = -9562.545399999999
